$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '30.288.05'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('E3').Value = '  -1.49%  '
Set-TextCell $ws.Range('D5') '237.91'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  +0.07%  '
Set-TextCell $ws.Range('D7') '0.4686'
$ws.Range('E7').Value = '  -1.93%  '
Set-TextCell $ws.Range('D8') '0.2831'
$ws.Range('E8').Value = '  -0.09%  '
Set-TextCell $ws.Range('D9') '0.06596'
$ws.Range('E9').Value = '  -1.60%  '
Set-TextCell $ws.Range('D10') '20.63'
$ws.Range('E10').Value = '  +10.06%  '
Set-TextCell $ws.Range('D11') '0.07785'
$ws.Range('E11').Value = '  +1.47%  '
Set-TextCell $ws.Range('D12') '98.12'
$ws.Range('E12').Value = '  -3.13%  '
Set-TextCell $ws.Range('D13') '1.883.78'
$ws.Range('E13').Value = '  -1.39%  '
Set-TextCell $ws.Range('D14') '5.091'
$ws.Range('E14').Value = '  -2.14%  '
Set-TextCell $ws.Range('D15') '0.6758'
$ws.Range('E15').Value = '  +1.06%  '
Set-TextCell $ws.Range('D16') '285.16'
$ws.Range('E16').Value = '  +9.97%  '
Set-TextCell $ws.Range('D17') '30.306.86'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('E18').Value = '  +0.02%  '
Set-TextCell $ws.Range('D19') '12.63'
$ws.Range('E19').Value = '  -0.09%  '
Set-TextCell $ws.Range('D20') '2.128.24'
$ws.Range('E20').Value = '  -1.37%  '
Set-TextCell $ws.Range('D21') '5.393'
$ws.Range('E21').Value = '  +0.08%  '
Set-TextCell $ws.Range('D22') '0.000007286'
$ws.Range('E22').Value = '  -2.41%  '
Set-TextCell $ws.Range('D23') '1.001'
$ws.Range('E23').Value = '  +0.10%  '
Set-TextCell $ws.Range('D24') '6.180'
$ws.Range('E24').Value = '  -1.59%  '
Set-TextCell $ws.Range('D25') '9.388'
$ws.Range('E25').Value = '  +0.39%  '
Set-TextCell $ws.Range('D26') '168.56'
$ws.Range('E26').Value = '  +0.70%  '
Set-TextCell $ws.Range('D27') '19.21'
$ws.Range('E27').Value = '  +0.40%  '
Set-TextCell $ws.Range('D28') '1.995'
$ws.Range('E28').Value = '  -3.03%  '
Set-TextCell $ws.Range('D29') '1.370'
$ws.Range('E29').Value = '  -1.11%  '
Set-TextCell $ws.Range('D30') '0.09711'
$ws.Range('E30').Value = '  -3.03%  '
Set-TextCell $ws.Range('D31') '4.402'
$ws.Range('E31').Value = '  -8.36%  '
$ws.Range('E32').Value = '  -1.58%  '
Set-TextCell $ws.Range('D33') '4.136'
$ws.Range('E33').Value = '  -2.81%  '
Set-TextCell $ws.Range('D34') '0.04666'
$ws.Range('E34').Value = '  -0.99%  '
Set-TextCell $ws.Range('D35') '0.7057'
$ws.Range('E35').Value = '  -2.65%  '
Set-TextCell $ws.Range('D36') '1.097'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws.Range('D37') '2.718'
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range('D38') '0.01872'
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range('D39') '6.639'
$ws.Range('E39').Value = '  +6.25%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Range('D40') '2.526'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range('D41') '72.06'
$ws.Range('E41').Value = '  -3.68%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range('D42') '1.973'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range('D43') '0.8663'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws.Range('D44') '1.000'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws.Range('D45') '103.00'
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws.Range('D46') '0.4188'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws.Range('D47') '987.09'
$ws.Range('E47').Value = '  +7.38%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range('D48') '7.276'
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range('D49') '9.231'
$ws.Range('E49').Value = '  +5.44%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell $ws.Range('D50') '34.02'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range('D51') '0.1146'
$ws.Range('E51').Value = '  -4.45%  '
